# myjobs, admin inf changes
#
# Applies 5 text replacements inside specific table cells.
#
# NOTE: this runtime's Range.Find.Execute(...) searches/replaces across the
# *entire* document regardless of which Range it is invoked on, so it is
# unsafe to use here (several of the target strings, e.g. "Hyderabad",
# "Not Applicable", "Campus Visit", repeat elsewhere in the document with
# different intended values). Instead we target each cell precisely via
# the Tables collection and assign to Range.Text directly, which correctly
# scopes the edit to that single cell and preserves the surrounding run
# formatting / any trailing empty runs.
#
#   Table 2 (Internship/Job details table):
#     r1c2  "Graphics Designer" -> "Software Developer"
#     r3c2  "Hyderabad"         -> "Delhi"
#     r5c2  "Hyderabad"         -> "Delhi"
#   Table 5 (IIT ISM Placement Calendar table):
#     r5c3 (Group Discussion row)   "Not Applicable" -> "Virtual"
#     r7c3 (Any other rounds row)   "Campus Visit"   -> "Not Applicable"

$d = $word.ActiveDocument

# --- Table 2: Job Designation / Place of Posting details ---
$t2 = $d.Tables.Item(2)

$t2.Cell(1, 2).Range.Text = "Software Developer"
$t2.Cell(3, 2).Range.Text = "Delhi"
$t2.Cell(5, 2).Range.Text = "Delhi"

# --- Table 5: IIT ISM Placement Calendar ---
$t5 = $d.Tables.Item(5)

$t5.Cell(5, 3).Range.Text = "Virtual"
$t5.Cell(7, 3).Range.Text = "Not Applicable"

Write-Output "Done applying replacements."
